$wb = $excel.ActiveWorkbook

# --- Create week 6 sheet by copying week 5's sheet (same template/layout) ---
$src = $wb.Worksheets.Item("Nädal 5")
$src.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count)) | Out-Null
$new = $wb.Worksheets.Item($wb.Worksheets.Count)
$new.Name = "Nädal 6"

# Fix the date-format style on B12 to match the other date cells in this sheet
# (the template row 12 used the generic style; week 6 needs an actual date there)
$new.Range("B8").Copy() | Out-Null
$new.Range("B12").PasteSpecial(-4122) | Out-Null

# --- Header: week-of date ---
$new.Range("G4").Value = 43535

# --- Row 7 ---
$new.Range("B7").Value = 43532
$new.Range("C7").Value = 0.33333333333333331
$new.Range("D7").Value = 0.43055555555555558
$new.Range("E7").Value = ""
$new.Range("G7").Value = "Class"
$new.Range("H7").Value = ""

# --- Row 8 (numeric/date/time cells first; text below, in authoring order) ---
$new.Range("B8").Value = 43533
$new.Range("C8").Value = 0.27083333333333331
$new.Range("D8").Value = 0.5
$new.Range("E8").Value = 60
$new.Range("G8").Value = "Prep."

# --- Row 9 ---
$new.Range("C9").Value = 0.76041666666666663
$new.Range("D9").Value = 0.78472222222222221

# New activity/comment text, entered in the order the author typed them
$new.Range("G9").Value = "Help"
$new.Range("H9").Value = "Helping other students"
$new.Range("H8").Value = "Teamtreehouse MVC course"

# --- Row 10 ---
$new.Range("B10").Value = 43534
$new.Range("C10").Value = 0.54166666666666663
$new.Range("D10").Value = 0.70138888888888884
$new.Range("E10").Value = ""
$new.Range("G10").Value = "Prep."
$new.Range("H10").Value = "Teamtreehouse React course"

# --- Row 11 ---
$new.Range("B11").Value = ""
$new.Range("C11").Value = 0.84027777777777779
$new.Range("D11").Value = 0.93055555555555547
$new.Range("G11").Value = "Prep."
$new.Range("H11").Value = "Working on my notes"
$new.Range("I11").Value = "x"

# --- Row 12 ---
$new.Range("B12").Value = 43535
$new.Range("C12").Value = 0.72916666666666663
$new.Range("D12").Value = 0.86805555555555547
$new.Range("G12").Value = "Help"
$new.Range("H12").Value = "Helping other students"

# Put the selection where the author left it and make week 6 the active tab
$new.Range("H13").Select() | Out-Null
$wb.Worksheets.Item("Nädal 6").Activate() | Out-Null
